$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet
#    (i.e. right after "2021-Q4"), holding the quarter's fund-holding
#    detail, mirroring the layout of the other quarterly sheets.
# -----------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the sibling quarterly sheets (0.75in
# left/right, 1in top/bottom, 0.5in header/footer -> expressed in points)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row (row 1)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row (row 2)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001703"
$newSheet.Range("C2").Value = "银华沪港深增长股票"
$newSheet.Range("D2").Value = "'3.02"
$newSheet.Range("E2").Value = "'87.52"
$newSheet.Range("F2").Value = "'5.11"
$newSheet.Range("G2").Value = "'0.1543"
$newSheet.Range("H2").Value = 3

# Header-style formatting: bold, centered, thin boxed border
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160
$newSheet.Range("B1:H1").Borders.LineStyle = 1

$newSheet.Range("A2").Font.Bold = $true
$newSheet.Range("A2").HorizontalAlignment = -4108
$newSheet.Range("A2").VerticalAlignment = -4160
$newSheet.Range("A2").Borders.LineStyle = 1

# -----------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add a new top data row for
#    2022-Q1 and shift the previous rows down, renumbering the
#    running index held in column A.
# -----------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.15

$totalSheet.Range("A2").Font.Bold = $true
$totalSheet.Range("A2").HorizontalAlignment = -4108
$totalSheet.Range("A2").VerticalAlignment = -4160
$totalSheet.Range("A2").Borders.LineStyle = 1

# Renumber the running index (column A) for the rows that were pushed down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
